$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

# 1) Duplicate the current last row (28, the blank "خلل" template row) down
#    to a new row 29, before filling row 28 in, so the new row keeps the
#    same blank template look (card=24, Correction=خلل, rest blank) that
#    row 28 used to have.
$ws.Rows.Item(28).Copy()
$ws.Rows.Item(29).Insert()

# 2) Now fill in row 28 (the new event being recorded), setting every
#    previously-empty cell to the text "nan" while leaving N28 ("خلل")
#    untouched.
$ws.Range("B28:M28").Value = "nan"
$ws.Range("O28:Q28").Value = "nan"
